$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) "We model Bretherton flow ..." paragraph (first bullet)
#    - insert "/Taylor" right after "Bretherton"
#    - change the trailing " by the diffuse interface model" to just " "
# ----------------------------------------------------------------------
$rng = $d.Range(0, $d.Content.End)
$rng.Find.Execute("Bretherton") | Out-Null
$insPoint = $d.Range($rng.End, $rng.End)
$insPoint.InsertAfter("/Taylor")

$d.Content.Find.Execute(
    " by the diffuse interface model", $false, $false, $false, $false, $false,
    $true, 1, $false, " ", 2) | Out-Null

# ----------------------------------------------------------------------
# 2) "Bubble can have axisymmetric ..." paragraph (second bullet)
#    - "... and non-axisymmetric shape depending ..." ->
#      "... and non-symmetric shape depending ..."
# ----------------------------------------------------------------------
$d.Content.Find.Execute(
    " and non-axisymmetric shape depending on the capillary number",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    " and non-symmetric shape depending on the capillary number", 2) | Out-Null

# ----------------------------------------------------------------------
# 3) "The lattice Boltzmann method (LBM) ..." paragraph (third bullet)
# ----------------------------------------------------------------------
$d.Content.Find.Execute(
    "The lattice Boltzmann method (LBM) free-energy binary liquid model was chosen as a framework for the diffuse interface model",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "We simulate with the lattice Boltzmann method (LBM) free-energy binary liquid model", 2) | Out-Null

# ----------------------------------------------------------------------
# 4) "Results show that ..." paragraph (fourth bullet)
# ----------------------------------------------------------------------
$d.Content.Find.Execute(
    "Results show that uniform density LBM free-energy binary liquid model is able to capture all phenomena indicated in literature as bubble shape change from non-axisymmetric shape to symmetric, existence of the vortex in front of a bubble, non-dimensional relative to liquid velocity of the bubble",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Results show that we capture all phenomena indicated in literature", 2) | Out-Null

# ----------------------------------------------------------------------
# 5) Replace the trailing empty paragraph with a new fifth bullet:
#    "Lattice Boltzmann method is a reliable tool to resolve dynamics of
#    Taylor bubbles "
#
#    InsertParagraphAfter(), called on the Range of the last real bullet,
#    duplicates that paragraph's pPr (ListParagraph style + numPr
#    ilvl=0/numId=2), so the freshly inserted paragraph keeps the correct
#    list numbering/style without minting a new numbering list. We put
#    the new sentence into that new paragraph, then delete the old empty
#    trailing paragraph that is no longer needed.
# ----------------------------------------------------------------------
$paraCountBefore = $d.Paragraphs.Count
$lastBulletPara = $d.Paragraphs.Item($paraCountBefore - 1)

$lastBulletPara.Range.InsertParagraphAfter()

$newBulletPara = $d.Paragraphs.Item($paraCountBefore)
$newInsPoint = $d.Range($newBulletPara.Range.Start, $newBulletPara.Range.Start)
$newInsPoint.InsertAfter("Lattice Boltzmann method is a reliable tool to resolve dynamics of Taylor bubbles ")

$newBulletPara = $d.Paragraphs.Item($paraCountBefore)
$oldTrailingEmptyPara = $d.Paragraphs.Item($paraCountBefore + 1)
$delRange = $d.Range($newBulletPara.Range.End - 1, $oldTrailingEmptyPara.Range.End)
$delRange.Delete()
